$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 50, pushing the
# existing rows 50-87 down to 52-89 (format is copied down from row 49,
# which already carries the date style used by column D).
$ws.Rows("50:51").Insert()

# --- New row 50: Artic Sprite / Primera -------------------------------
$ws.Range("A50").Value = 1
$ws.Range("B50").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C50").Value = "Arica y Parinacota"
$ws.Range("D50").Value = 44960
$ws.Range("E50").Value = 15
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100103
$ws.Range("H50").Value = "Frutos de hueso (carozo)"
$ws.Range("I50").Value = 100103006
$ws.Range("J50").Value = "Nectarín"
$ws.Range("K50").Value = "Artic Sprite"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 300
$ws.Range("N50").Value = 24000
$ws.Range("O50").Value = 25000
$ws.Range("P50").Value = 24500
$ws.Range("Q50").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 1361
$ws.Range("T50").Value = 18

# --- New row 51: Venus / Primera --------------------------------------
$ws.Range("A51").Value = 1
$ws.Range("B51").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C51").Value = "Arica y Parinacota"
$ws.Range("D51").Value = 44960
$ws.Range("E51").Value = 15
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100103
$ws.Range("H51").Value = "Frutos de hueso (carozo)"
$ws.Range("I51").Value = 100103006
$ws.Range("J51").Value = "Nectarín"
$ws.Range("K51").Value = "Venus"
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 250
$ws.Range("N51").Value = 24000
$ws.Range("O51").Value = 25000
$ws.Range("P51").Value = 24500
$ws.Range("Q51").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R51").Value = "Región de O'Higgins"
$ws.Range("S51").Value = 1361
$ws.Range("T51").Value = 18
